$d = $word.ActiveDocument

# Helper: rename an InlineShape's drawing object name (wp:docPr/@name).
# InlineShape itself exposes no writable Name in the Word object model,
# so the picture is briefly promoted to a floating Shape (which does
# expose .Name), renamed, then converted back to an inline picture in
# place.
function Rename-InlineShape($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

$section = $d.Sections.Item(1)

# First-page footer (footer1.xml): Pearson logo "image2.png" -> "image1.png"
$firstFooter = $section.Footers.Item(2)
Rename-InlineShape $firstFooter.Range.InlineShapes.Item(1) "image1.png"

# Default footer (footer2.xml): Pearson logo "image2.png" -> "image1.png"
$defaultFooter = $section.Footers.Item(1)
Rename-InlineShape $defaultFooter.Range.InlineShapes.Item(1) "image1.png"

# First-page header (header1.xml): BTEC logo "image1.jpg" -> "image2.jpg"
$firstHeader = $section.Headers.Item(2)
Rename-InlineShape $firstHeader.Range.InlineShapes.Item(1) "image2.jpg"
